# Append a new "05-aug" column (BA) to the "Prix Spot" sheet, mirroring the
# existing AZ ("04-aug") column's layout: a bold/bordered/centered header in
# row 1 and numeric values in rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Header cell - reuse the formatting of the preceding date header (AZ1),
# which carries the bold/centered/bordered header style.
$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)
$ws.Range("BA1").Value = "05-aug"

# Hourly values for 05-aug, row by row (row 2 = "00 - 01" ... row 25 = "23 - 24").
$values = @(
    8.460000000000001,
    -0.07000000000000001,
    -0.01,
    -0.06,
    -0.25,
    -0.11,
    0.02,
    2.77,
    0,
    -0.03,
    -2.78,
    -3.54,
    -9.65,
    -24.02,
    -19.5,
    -14.94,
    -14.44,
    -0.02,
    3,
    46.4,
    73.5,
    82.56,
    88.22,
    74.95
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 53).Value = $values[$i]
}
